$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194, shifting existing rows 194:280 down to 195:281
$ws.Rows.Item(194).Insert()

# Populate the new row 194 with the latest weekly data point.
# Columns A, B, C, E, F, G, H, I, O, R are constant for this dataset.
$ws.Cells.Item(194, 1).Value = 8
$ws.Cells.Item(194, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(194, 3).Value = "Coquimbo"
$ws.Cells.Item(194, 4).Value = 45141
$ws.Cells.Item(194, 5).Value = 4
$ws.Cells.Item(194, 6).Value = 100112001
$ws.Cells.Item(194, 7).Value = "Berenjena"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 440
$ws.Cells.Item(194, 11).Value = 9000
$ws.Cells.Item(194, 12).Value = 10000
$ws.Cells.Item(194, 13).Value = 9500
$ws.Cells.Item(194, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(194, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(194, 16).Value = 190
$ws.Cells.Item(194, 17).Value = 50
$ws.Cells.Item(194, 18).Value = "Hortaliza"
